$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-30 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-01 Wednesday", 2)

$d.Content.Find.Execute("798×3=", $true, $false, $false, $false, $false, $true, 1, $false, "265×4=", 2)
$d.Content.Find.Execute("336×3=", $true, $false, $false, $false, $false, $true, 1, $false, "114×7=", 2)
$d.Content.Find.Execute("104×2=", $true, $false, $false, $false, $false, $true, 1, $false, "245×9=", 2)
$d.Content.Find.Execute("607×4=", $true, $false, $false, $false, $false, $true, 1, $false, "643×4=", 2)
$d.Content.Find.Execute("566×3=", $true, $false, $false, $false, $false, $true, 1, $false, "495×4=", 2)

$d.Content.Find.Execute("606×5=", $true, $false, $false, $false, $false, $true, 1, $false, "307×7=", 2)
$d.Content.Find.Execute("374×3=", $true, $false, $false, $false, $false, $true, 1, $false, "979×9=", 2)
$d.Content.Find.Execute("170×7=", $true, $false, $false, $false, $false, $true, 1, $false, "648×4=", 2)
$d.Content.Find.Execute("840×4=", $true, $false, $false, $false, $false, $true, 1, $false, "734×5=", 2)
$d.Content.Find.Execute("431×7=", $true, $false, $false, $false, $false, $true, 1, $false, "826×2=", 2)

$d.Content.Find.Execute("402×3=", $true, $false, $false, $false, $false, $true, 1, $false, "836×4=", 2)
$d.Content.Find.Execute("455×9=", $true, $false, $false, $false, $false, $true, 1, $false, "261×4=", 2)
$d.Content.Find.Execute("771×4=", $true, $false, $false, $false, $false, $true, 1, $false, "303×6=", 2)
$d.Content.Find.Execute("911×7=", $true, $false, $false, $false, $false, $true, 1, $false, "694×2=", 2)
$d.Content.Find.Execute("496×4=", $true, $false, $false, $false, $false, $true, 1, $false, "139×8=", 2)

$d.Content.Find.Execute("655×8=", $true, $false, $false, $false, $false, $true, 1, $false, "434×4=", 2)
$d.Content.Find.Execute("724×8=", $true, $false, $false, $false, $false, $true, 1, $false, "580×5=", 2)
$d.Content.Find.Execute("538×2=", $true, $false, $false, $false, $false, $true, 1, $false, "607×5=", 2)
$d.Content.Find.Execute("262×9=", $true, $false, $false, $false, $false, $true, 1, $false, "424×3=", 2)
$d.Content.Find.Execute("449×7=", $true, $false, $false, $false, $false, $true, 1, $false, "700×9=", 2)

$d.Content.Find.Execute("708×5=", $true, $false, $false, $false, $false, $true, 1, $false, "160×5=", 2)
$d.Content.Find.Execute("691×5=", $true, $false, $false, $false, $false, $true, 1, $false, "638×4=", 2)
$d.Content.Find.Execute("912×9=", $true, $false, $false, $false, $false, $true, 1, $false, "272×7=", 2)
$d.Content.Find.Execute("669×3=", $true, $false, $false, $false, $false, $true, 1, $false, "316×8=", 2)
$d.Content.Find.Execute("879×4=", $true, $false, $false, $false, $false, $true, 1, $false, "926×8=", 2)
